$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (test case #1 header row) ---
$ws.Range("C2").Value = "no input"
# D2/E2 unchanged ("using initialization list")

# --- Row 3: becomes "5 d" / "5 d" / PASS, with a time-style number format on D3:E3 ---
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "5 d"
$ws.Range("E3").Value = "5 d"
$ws.Range("D3:E3").NumberFormat = "h:mm AM/PM"
$ws.Range("F3").Value = "PASS"

# --- Row 4: becomes "10 f" / "10 f" ---
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = "10 f"
$ws.Range("E4").Value = "10 f"

# --- Old rows 5-9 content is removed entirely ---
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()

# --- Old row 10 content is replaced by a brand-new test case (row 10) ---
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "to initialize private variables using constructors"
$ws.Range("C10").Value = "no input"
$ws.Range("D10").Value = "100 d"
$ws.Range("E10").Value = "100 d"
$ws.Range("F10").Value = "PASS"

# --- New row 11 ---
$ws.Range("D11").Value = "10 f"
$ws.Range("E11").Value = "10 f"

# --- Move the active selection to F1 ---
$ws.Range("F1").Select()
